# "update tasks for Tet Holiday"
# Sheet1 holds the task tracker. Three tasks (rows 22, 27, 28) that had no
# owner and were still marked "waiting" get assigned to Hung and flipped to
# "on processing" now that work resumes after the Tet break.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 22 - "Bat loi xuat hien trong qua trinh xu ly cau hoi"
$ws.Range("C22").Value = "Hung"
$ws.Range("G22").Value = "on processing"

# Row 27 - "Trinh bay noi dung hien thi trong tim kiem bang hoi dap"
$ws.Range("C27").Value = "Hung"
$ws.Range("G27").Value = "on processing"

# Row 28 - "Viet Hoa QAWeb"
$ws.Range("C28").Value = "Hung"
$ws.Range("G28").Value = "on processing"

# Scroll the sheet view down and move the active selection, matching the
# reviewer's position when they made this pass over the sheet.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C35").Select()
